# This edit rearranges the weekly records in rows 2-14 (columns D, J, K, L, M, O, P)
# so that each row's price/volume/date/origin data moves to a different row.
# Row 12 is untouched. We capture all original values first (using Value2() for
# the date column D so we keep the raw serial number instead of a DateTime, and
# Value() for everything else), then write them back out in their new positions
# so the whole operation is safe regardless of write order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping: destination row -> source row (i.e. the data that used to live in
# the source row now lives in the destination row).
$mapping = @{
    2  = 4
    3  = 5
    4  = 14
    5  = 6
    6  = 10
    7  = 13
    8  = 3
    9  = 11
    10 = 2
    11 = 9
    12 = 12
    13 = 7
    14 = 8
}

# Columns that carry per-row data which gets shuffled: D, J, K, L, M, O, P
$cols = @(4, 10, 11, 12, 13, 15, 16)

# Snapshot original values for every row/column involved before writing anything.
$original = @{}
foreach ($row in 2..14) {
    $rowData = @{}
    foreach ($col in $cols) {
        if ($col -eq 4) {
            # Date column: read the raw serial number, not a converted DateTime.
            $rowData[$col] = $ws.Cells.Item($row, $col).Value2()
        } else {
            $rowData[$col] = $ws.Cells.Item($row, $col).Value()
        }
    }
    $original[$row] = $rowData
}

# Write the values to their new destination rows using the captured snapshot.
foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $srcData = $original[$srcRow]
    foreach ($col in $cols) {
        $ws.Cells.Item($destRow, $col).Value = $srcData[$col]
    }
}
